$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 keeps showing "Meetings" (unchanged text) and B4 keeps showing the
# "manage the money" task, just reworded to "Manage money" right before
# it is turned into a hyperlink.
$ws.Range("B3").Value = "Meetings"
$ws.Range("B4").Value = "Manage money"

# Turn B4 into a hyperlink. This also registers Excel's built-in
# "Hyperlink" cell style (underlined, theme-colored font) and applies it
# to B4, matching the new font/style entries added to styles.xml.
$ws.Hyperlinks.Add($ws.Range("B4"), "https://example.com/")

[void]$ws.Range("B4").Select()
